$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rscKpi code labels: GP1/GP2/GP3 -> GP01/GP02/GP03 and BP1/BP2/BP3 -> BP01/BP02/BP03
# These values live in column B, rows 16-60.
$ws.Range("B16").Value = "GP01"
$ws.Range("B17").Value = "GP02"
$ws.Range("B18").Value = "GP02"
$ws.Range("B19").Value = "GP03"
$ws.Range("B20").Value = "GP03"

$ws.Range("B21:B25").Value = "BP01"
$ws.Range("B26:B30").Value = "BP02"
$ws.Range("B31:B60").Value = "BP03"

# Reset the sheet view back to the top-left corner (A1) instead of G1/I4.
$ws.Range("A1").Select()
